# Towers vs. Enemies chart.xlsx - expand the tower/enemy stat tables.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Make room: push the old "Enemies Money Dropped" mini table (rows 33-45)
#    down to rows 48-60, and open one extra row just below row 21 so the
#    "Cost:" table grows from 10 rows to 11 rows (blank row + header + 9 rows).
# ---------------------------------------------------------------------------
$ws.Range("A22:A22").EntireRow.Insert()
$ws.Range("A34:A47").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) Row 22 is now a blank spacer row (style already copied down from row 21
#    by the insert above, matching the A/B/C/J/K/L quote-prefix styling).
#    Clear any stray value it might carry.
# ---------------------------------------------------------------------------
$ws.Range("A22:L22").ClearContents()

# ---------------------------------------------------------------------------
# 3) Row 23: new header row for the expanded tower table.
# ---------------------------------------------------------------------------
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "Cost:"
$ws.Range("D23").Value = "Health"
$ws.Range("F23").Value = "Armor damage"
$ws.Range("G23").Value = "projectile Speed"
$ws.Range("C23").Value = "Damage (per sec)"
$ws.Range("H23").Value = "Fire Rate"
$ws.Range("E23").Value = "Range(radius)"
$ws.Range("I23").Value = "Money Increase(per sec)"
$ws.Range("J23").Value = "Happieness"

# ---------------------------------------------------------------------------
# 4) Rows 24-32: tower data (name, cost, dmg/sec, health, range, armor dmg,
#    projectile speed, fire rate, money increase/sec, happiness).
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "House"
$ws.Range("B24").Value = 300
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = "fast"
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 3

$ws.Range("A25").Value = "Recycling Centers"
$ws.Range("B25").Value = 450
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 150
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = "med"
$ws.Range("I25").Value = -2
$ws.Range("J25").Value = 0

$ws.Range("A26").Value = "Water Purifacation"
$ws.Range("B26").Value = 450
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 150
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = "med"
$ws.Range("I26").Value = -2
$ws.Range("J26").Value = 0

$ws.Range("A27").Value = "GreenBelt/Park"
$ws.Range("B27").Value = 500
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 150
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = "med"
$ws.Range("I27").Value = -1
$ws.Range("J27").Value = 7

$ws.Range("A28").Value = "Schools/Youth Centers"
$ws.Range("B28").Value = 700
$ws.Range("C28").Value = 15
$ws.Range("D28").Value = 250
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = "fast"
$ws.Range("I28").Value = -3
$ws.Range("J28").Value = 2

$ws.Range("A29").Value = "Police/Fire Department"
$ws.Range("B29").Value = 1500
$ws.Range("C29").Value = 35
$ws.Range("D29").Value = 500
$ws.Range("E29").Value = 2
$ws.Range("F29").Value = 3
$ws.Range("G29").NumberFormat = "d-mmm"
$ws.Range("G29").Value = "One Half"
$ws.Range("H29").Value = "slow"
$ws.Range("I29").Value = -2
$ws.Range("J29").Value = 1

$ws.Range("A30").Value = "Stores"
$ws.Range("B30").Value = 650
$ws.Range("C30").Value = 10
$ws.Range("D30").Value = 200
$ws.Range("E30").Value = 3
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = "fast"
$ws.Range("I30").Value = 10
$ws.Range("J30").Value = 2

$ws.Range("A31").Value = "Industry"
$ws.Range("B31").Value = 550
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 200
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = "slow"
$ws.Range("I31").Value = 8
$ws.Range("J31").Value = -10

$ws.Range("A32").Value = "National Parks/Monuments"
$ws.Range("B32").Value = 2000
$ws.Range("C32").Value = "1 or 0"
$ws.Range("D32").Value = 1000
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 10
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = "slow"
$ws.Range("I32").Value = 35
$ws.Range("J32").Value = 7

# ---------------------------------------------------------------------------
# 5) Row 34: header for the new enemy Speed/Damage/Health/Armor mini table.
# ---------------------------------------------------------------------------
$ws.Range("B34").Value = "Speed"
$ws.Range("E34").Value = "Armor"
$ws.Range("C34").Value = "Damage(per sec)"
$ws.Range("D34").Value = "Health"

# ---------------------------------------------------------------------------
# 6) Rows 35-45: enemy data (name, speed, dmg/sec, health, armor).
# ---------------------------------------------------------------------------
$ws.Range("A35").Value = "Smog"
$ws.Range("B35").Value = 2
$ws.Range("C35").Value = 3
$ws.Range("D35").Value = 35
$ws.Range("E35").Value = 0

$ws.Range("A36").Value = "Gangs"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 5
$ws.Range("D36").Value = 50
$ws.Range("E36").Value = 20

$ws.Range("A37").Value = "Arsonist"
$ws.Range("B37").Value = 3
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 20
$ws.Range("E37").Value = 5

$ws.Range("A38").Value = "Joe-Bloe Criminal"
$ws.Range("B38").Value = 3
$ws.Range("C38").Value = 3
$ws.Range("D38").Value = 25
$ws.Range("E38").Value = 5

$ws.Range("A39").Value = "Trash Buildup"
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 7
$ws.Range("D39").Value = 60
$ws.Range("E39").Value = 10

$ws.Range("A40").Value = "Water Polution"
$ws.Range("B40").Value = 2
$ws.Range("C40").Value = 5
$ws.Range("D40").Value = 35
$ws.Range("E40").Value = 0

$ws.Range("A41").Value = "Fire"
$ws.Range("B41").Value = "3 or 4"
$ws.Range("C41").Value = 2
$ws.Range("D41").Value = 20
$ws.Range("E41").Value = 0

$ws.Range("A42").Value = "Flood"
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 35
$ws.Range("E42").Value = 0

$ws.Range("A43").Value = "Earthquake"
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = 4
$ws.Range("D43").Value = 45
$ws.Range("E43").Value = 0

$ws.Range("A44").Value = "Graffiti"
$ws.Range("B44").Value = 2
$ws.Range("C44").Value = 3
$ws.Range("D44").Value = 35
$ws.Range("E44").Value = 5

$ws.Range("A45").Value = "Poor Education"
$ws.Range("B45").Value = 2
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 30
$ws.Range("E45").Value = 5

# ---------------------------------------------------------------------------
# 7) Column widths for the newly-used columns H, I, K (J already existed).
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 17.7109375
$ws.Columns.Item(9).ColumnWidth = 17.8671875
$ws.Columns.Item(11).ColumnWidth = 18.1666666666667

# ---------------------------------------------------------------------------
# 8) View state: scroll position + active selection, matching the saved file.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E46").Select()
